# Auto-generated Excel COM-interop script
# Applies updated market-board figures (currentAveragePrice / LevePrice / LeveProfit columns)
# to the Midgardsormr_Profits workbook, per sheet (crafting job) and leve row.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2640
$ws.Range("I40").Value = 1800
$ws.Range("K40").Value = 1800
$ws.Range("M40").Value = -1625
$ws.Range("H112").Value = 6770.766
$ws.Range("I112").Value = 1105.5
$ws.Range("J112").Value = 7022.5557
$ws.Range("K112").Value = 3316.5
$ws.Range("L112").Value = 21067.6671
$ws.Range("M112").Value = -2208.5
$ws.Range("N112").Value = -23283.6671
$ws.Range("H133").Value = 73333
$ws.Range("H138").Value = 3197.2593
$ws.Range("I138").Value = 2647.9375
$ws.Range("J138").Value = 3996.2727
$ws.Range("K138").Value = 7943.8125
$ws.Range("L138").Value = 11988.8181
$ws.Range("M138").Value = -2803.8125
$ws.Range("N138").Value = -22268.8181

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 6001
$ws.Range("J8").Value = 8000
$ws.Range("L8").Value = 8000
$ws.Range("N8").Value = -8288
$ws.Range("H32").Value = 17715.012
$ws.Range("I32").Value = 17557.514
$ws.Range("K32").Value = 17557.514
$ws.Range("M32").Value = -17270.514
$ws.Range("H132").Value = 3361.75
$ws.Range("I132").Value = 3146.6667
$ws.Range("K132").Value = 9440.000100000001
$ws.Range("M132").Value = -6910.000100000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("L35").ClearContents()
$ws.Range("H86").Value = 1276.7084
$ws.Range("I86").Value = 1237.35
$ws.Range("J86").Value = 1473.5
$ws.Range("K86").Value = 1237.35
$ws.Range("L86").Value = 1473.5
$ws.Range("M86").Value = -114.3499999999999
$ws.Range("N86").Value = -3719.5
$ws.Range("H89").Value = 1276.7084
$ws.Range("I89").Value = 1237.35
$ws.Range("J89").Value = 1473.5
$ws.Range("K89").Value = 6186.75
$ws.Range("L89").Value = 7367.5
$ws.Range("M89").Value = -570.75
$ws.Range("N89").Value = -18599.5
$ws.Range("H105").Value = 3134.2307
$ws.Range("I105").Value = 1305.1111
$ws.Range("J105").Value = 7249.75
$ws.Range("K105").Value = 1305.1111
$ws.Range("L105").Value = 7249.75
$ws.Range("M105").Value = 441.8888999999999
$ws.Range("N105").Value = -10743.75
$ws.Range("H134").Value = 4782.5884
$ws.Range("I134").Value = 4850.357
$ws.Range("K134").Value = 14551.071
$ws.Range("M134").Value = -12016.071

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 10139.167
$ws.Range("I62").Value = 4811.75
$ws.Range("K62").Value = 4811.75
$ws.Range("M62").Value = -4187.75
$ws.Range("H65").Value = 10139.167
$ws.Range("I65").Value = 4811.75
$ws.Range("K65").Value = 24058.75
$ws.Range("M65").Value = -20938.75
$ws.Range("H86").Value = 37767.09
$ws.Range("I86").Value = 76608
$ws.Range("K86").Value = 76608
$ws.Range("M86").Value = -75485
$ws.Range("H89").Value = 37767.09
$ws.Range("I89").Value = 76608
$ws.Range("K89").Value = 383040
$ws.Range("M89").Value = -377424
$ws.Range("H106").Value = 10000
$ws.Range("I106").Value = 10000
$ws.Range("K106").Value = 10000
$ws.Range("M106").Value = -8738
$ws.Range("H107").Value = 343.03845
$ws.Range("I107").Value = 279.1875
$ws.Range("J107").Value = 445.2
$ws.Range("K107").Value = 279.1875
$ws.Range("L107").Value = 445.2
$ws.Range("M107").Value = 1640.8125
$ws.Range("N107").Value = -4285.2
$ws.Range("H132").Value = 48156
$ws.Range("I132").Value = 56224.816
$ws.Range("K132").Value = 168674.448
$ws.Range("M132").Value = -166144.448
$ws.Range("H138").Value = 118998.664
$ws.Range("J138").Value = 118998.664
$ws.Range("L138").Value = 118998.664
$ws.Range("N138").Value = -129278.664

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H133").Value = 6328.625
$ws.Range("J133").Value = 9900
$ws.Range("L133").Value = 29700
$ws.Range("N133").Value = -39820

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 22000
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H20").Value = 46221.777
$ws.Range("J20").Value = 49999.5
$ws.Range("L20").Value = 49999.5
$ws.Range("N20").Value = -50489.5
$ws.Range("H24").Value = 21525.736
$ws.Range("J24").Value = 23999.375
$ws.Range("L24").Value = 23999.375
$ws.Range("N24").Value = -24345.375
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("N88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("N91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("H126").Value = 2695.7827
$ws.Range("I126").Value = 2587.842
$ws.Range("K126").Value = 7763.526
$ws.Range("M126").Value = -5293.526
$ws.Range("H132").Value = 1436.8462
$ws.Range("I132").Value = 1075.3
$ws.Range("J132").Value = 2642
$ws.Range("K132").Value = 3225.9
$ws.Range("L132").Value = 7926
$ws.Range("M132").Value = -695.8999999999996
$ws.Range("N132").Value = -12986

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 11868.77
$ws.Range("J122").Value = 13968.8
$ws.Range("L122").Value = 41906.39999999999
$ws.Range("N122").Value = -46806.39999999999
$ws.Range("H132").Value = 2420.2
$ws.Range("I132").Value = 2521.4348
$ws.Range("K132").Value = 7564.3044
$ws.Range("M132").Value = -5034.3044
$ws.Range("H136").Value = 3094.7036
$ws.Range("I136").Value = 2873.9285
$ws.Range("K136").Value = 8621.7855
$ws.Range("M136").Value = -6071.7855

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1467.25
$ws.Range("I107").Value = 1934.5
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 5803.5
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -3883.5
$ws.Range("N107").Value = -6840
$ws.Range("H122").Value = 49128.105
$ws.Range("I122").Value = 102188.695
$ws.Range("K122").Value = 306566.085
$ws.Range("M122").Value = -304116.085
$ws.Range("H126").Value = 201568
$ws.Range("I126").Value = 1509.5714
$ws.Range("K126").Value = 4528.7142
$ws.Range("M126").Value = -2058.7142
$ws.Range("H132").Value = 18817.377
$ws.Range("I132").Value = 23651.629
$ws.Range("K132").Value = 70954.887
$ws.Range("M132").Value = -68424.887
$ws.Range("H136").Value = 25409.812
$ws.Range("I136").Value = 27859.861
$ws.Range("K136").Value = 83579.583
$ws.Range("M136").Value = -81029.583
